# Update test artefacts to new format of fiscal years
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1 headers: rename the fiscal-year / quarter columns to the new "FY20xx" format
$ws.Range("T1").Value = "FY2019 Q1 (D)"
$ws.Range("U1").Value = "FY2018 Q3 (D)"
$ws.Range("V1").Value = "FY2018 Q2 (D)"
$ws.Range("W1").Value = "FY2018 Q1 (D)"
$ws.Range("X1").Value = "FY2018 (MTEF)"

# Update the view: scroll so column O is the left-most visible column, and
# move the active selection to Z5:Z6
$excel.ActiveWindow.ScrollColumn = 15
$ws.Range("Z5:Z6").Select()
